# Appends new execution-log rows to each TestCase sheet (name path changed as wait not working)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TestCase1_HospitalFilter")
$ws1.Range("A161").Value = 'Opened Practo homepage.'
$ws1.Range("A162").Value = 'Searching for hospitals in: Bangalore'
$ws1.Range("A163").Value = 'Error during search: Expected condition failed: waiting for element found by By.xpath: //div[contains(text(),''Bangalore'')] to be clickable, but the element was not found: org.openqa.selenium.NoSuchElementException: no such element: Unable to locate element: {"method":"xpath","selector":"//div[contains(text(),''Bangalore'')]"}.
(tried for 15 seconds with 500 milliseconds interval)
Build info: version: ''4.40.0'', revision: ''b3333f1''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''21''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 142.0.7444.176, chrome: {chromedriverVersion: 142.0.7444.175 (302067f14a4..., userDataDir: C:\Users\2457382\AppData\Lo...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:58569}, goog:processID: 11444, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:58569/devtoo..., se:cdpVersion: 142.0.7444.176, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 86d10fbe4c0ea41f2bc6d58ba35182d5'
$ws1.Range("A164").Value = 'Printing and Writing'
$ws1.Range("A165").Value = 'Opened Practo homepage.'
$ws1.Range("A166").Value = 'Searching for hospitals in: Bangalore'
$ws1.Range("A167").Value = 'Opened Practo homepage.'
$ws1.Range("A168").Value = 'Searching for hospitals in: Bangalore'
$ws1.Range("A169").Value = 'Error during search: invalid session id: session deleted as the browser has closed the connection
from disconnected: not connected to DevTools
  (Session info: chrome=142.0.7444.176)
Build info: version: ''4.40.0'', revision: ''b3333f1''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''21''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [28dc2a14f3dabfc3d7ac1b4d2d4b0598, sendKeysToElement {id=f.0EC9F77330D5E8046A6531791E6535AA.d.18207EDF1F27767DC78E19B2085C4269.e.2, value=[Ljava.lang.CharSequence;@1934339}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 142.0.7444.176, chrome: {chromedriverVersion: 142.0.7444.175 (302067f14a4..., userDataDir: C:\Users\2457382\AppData\Lo...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:57189}, goog:processID: 3700, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:57189/devtoo..., se:cdpVersion: 142.0.7444.176, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Element: [[ChromeDriver: chrome on windows (28dc2a14f3dabfc3d7ac1b4d2d4b0598)] -> xpath: //input[@placeholder=''Search location'']]
Session ID: 28dc2a14f3dabfc3d7ac1b4d2d4b0598'
$ws1.Range("A170").Value = 'Scrolling error: invalid session id
Build info: version: ''4.40.0'', revision: ''b3333f1''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''21''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [28dc2a14f3dabfc3d7ac1b4d2d4b0598, executeScript {script=window.scrollBy(0, 1000), args=[]}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 142.0.7444.176, chrome: {chromedriverVersion: 142.0.7444.175 (302067f14a4..., userDataDir: C:\Users\2457382\AppData\Lo...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:57189}, goog:processID: 3700, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:57189/devtoo..., se:cdpVersion: 142.0.7444.176, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 28dc2a14f3dabfc3d7ac1b4d2d4b0598'
$ws1.Range("A171").Value = 'Opened Practo homepage.'
$ws1.Range("A172").Value = 'Searching for hospitals in: Bangalore'
$ws1.Range("A173").Value = 'Added hospital: Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws1.Range("A174").Value = 'Added hospital: Manipal Hospitals'
$ws1.Range("A175").Value = 'Added hospital: Koshys Hospital'
$ws1.Range("A176").Value = 'Added hospital: Motherhood Hospital'
$ws1.Range("A177").Value = 'Added hospital: Motherhood Hospital'
$ws1.Range("A178").Value = 'Added hospital: Trilife Hospital'
$ws1.Range("A179").Value = 'Added hospital: Apollo Cradle & Children’s Hospital'
$ws1.Range("A180").Value = 'Printing and Writing'
$ws1.Range("A181").Value = 'Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws1.Range("A182").Value = 'Manipal Hospitals'
$ws1.Range("A183").Value = 'Koshys Hospital'
$ws1.Range("A184").Value = 'Motherhood Hospital'
$ws1.Range("A185").Value = 'Motherhood Hospital'
$ws1.Range("A186").Value = 'Trilife Hospital'
$ws1.Range("A187").Value = 'Apollo Cradle & Children’s Hospital'
$ws1.Range("A188").Value = 'Opened Practo homepage.'
$ws1.Range("A189").Value = 'Searching for hospitals in: Bangalore'
$ws1.Range("A190").Value = 'Added hospital: Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws1.Range("A191").Value = 'Added hospital: Manipal Hospitals'
$ws1.Range("A192").Value = 'Added hospital: Koshys Hospital'
$ws1.Range("A193").Value = 'Added hospital: Motherhood Hospital'
$ws1.Range("A194").Value = 'Added hospital: Motherhood Hospital'
$ws1.Range("A195").Value = 'Added hospital: Trilife Hospital'
$ws1.Range("A196").Value = 'Added hospital: Apollo Cradle & Children’s Hospital'
$ws1.Range("A197").Value = 'Printing and Writing'
$ws1.Range("A198").Value = 'Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws1.Range("A199").Value = 'Manipal Hospitals'
$ws1.Range("A200").Value = 'Koshys Hospital'
$ws1.Range("A201").Value = 'Motherhood Hospital'
$ws1.Range("A202").Value = 'Motherhood Hospital'
$ws1.Range("A203").Value = 'Trilife Hospital'
$ws1.Range("A204").Value = 'Apollo Cradle & Children’s Hospital'

$ws2 = $wb.Worksheets.Item("TestCase0_MaxRatingFinder")
$ws2.Range("A57").Value = 'Searching for hospitals in: Bangalore'
$ws2.Range("A58").Value = 'Checking for rating ≥ 4.50: Bangalore'
$ws2.Range("A59").Value = 'Rating captured: 4.50'
$ws2.Range("A60").Value = 'PASS — Rating ≥ 4.50 (actual: 4.50)'
$ws2.Range("A61").Value = 'Result: Max Rating meets threshold (≥ 4.50): actual 4.50'
$ws2.Range("A62").Value = 'Searching for hospitals in: Bangalore'
$ws2.Range("A63").Value = 'Checking for rating ≥ 4.50: Bangalore'
$ws2.Range("A64").Value = 'Rating captured: 4.50'
$ws2.Range("A65").Value = 'PASS — Rating ≥ 4.50 (actual: 4.50)'
$ws2.Range("A66").Value = 'Result: Max Rating meets threshold (≥ 4.50): actual 4.50'

$ws3 = $wb.Worksheets.Item("TestCase2_TopCities")
$ws3.Range("A100").Value = 'Top Cities:'
$ws3.Range("A101").Value = 'Bangalore'
$ws3.Range("A102").Value = 'Delhi'
$ws3.Range("A103").Value = 'Mumbai'
$ws3.Range("A104").Value = 'Chennai'
$ws3.Range("A105").Value = 'Hyderabad'
$ws3.Range("A106").Value = 'Kolkata'
$ws3.Range("A107").Value = 'Pune'
$ws3.Range("A108").Value = 'Ahmedabad'
$ws3.Range("A109").Value = 'Top Cities:'
$ws3.Range("A110").Value = 'Bangalore'
$ws3.Range("A111").Value = 'Delhi'
$ws3.Range("A112").Value = 'Mumbai'
$ws3.Range("A113").Value = 'Chennai'
$ws3.Range("A114").Value = 'Hyderabad'
$ws3.Range("A115").Value = 'Kolkata'
$ws3.Range("A116").Value = 'Pune'
$ws3.Range("A117").Value = 'Ahmedabad'
$ws3.Range("A118").Value = 'Top Cities:'
$ws3.Range("A119").Value = 'Bangalore'
$ws3.Range("A120").Value = 'Delhi'
$ws3.Range("A121").Value = 'Mumbai'
$ws3.Range("A122").Value = 'Chennai'
$ws3.Range("A123").Value = 'Hyderabad'
$ws3.Range("A124").Value = 'Kolkata'
$ws3.Range("A125").Value = 'Pune'
$ws3.Range("A126").Value = 'Ahmedabad'
$ws3.Range("A127").Value = 'Top Cities:'
$ws3.Range("A128").Value = 'Bangalore'
$ws3.Range("A129").Value = 'Delhi'
$ws3.Range("A130").Value = 'Mumbai'
$ws3.Range("A131").Value = 'Chennai'
$ws3.Range("A132").Value = 'Hyderabad'
$ws3.Range("A133").Value = 'Kolkata'
$ws3.Range("A134").Value = 'Pune'
$ws3.Range("A135").Value = 'Ahmedabad'

$ws4 = $wb.Worksheets.Item("TestCase3_InvalidForm")
$ws4.Range("A84").Value = 'Forced click on Schedule button.'
$ws4.Range("A85").Value = 'No error messages found.'
$ws4.Range("A86").Value = 'Checking for invalid fields...'
$ws4.Range("A87").Value = 'Empty Name'
$ws4.Range("A88").Value = 'Empty Organization Name'
$ws4.Range("A89").Value = 'Invalid Contact Number'
$ws4.Range("A90").Value = 'Invalid Email ID'
$ws4.Range("A91").Value = 'Forced click on Schedule button.'
$ws4.Range("A92").Value = 'No error messages found.'
$ws4.Range("A93").Value = 'Checking for invalid fields...'
$ws4.Range("A94").Value = 'Empty Name'
$ws4.Range("A95").Value = 'Empty Organization Name'
$ws4.Range("A96").Value = 'Invalid Contact Number'
$ws4.Range("A97").Value = 'Invalid Email ID'
$ws4.Range("A98").Value = 'Forced click on Schedule button.'
$ws4.Range("A99").Value = 'No error messages found.'
$ws4.Range("A100").Value = 'Checking for invalid fields...'
$ws4.Range("A101").Value = 'Empty Name'
$ws4.Range("A102").Value = 'Empty Organization Name'
$ws4.Range("A103").Value = 'Invalid Contact Number'
$ws4.Range("A104").Value = 'Invalid Email ID'

